# Fill in two new rows (37 and 38) of coding-question tracker data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: "Maximum circular subarray sum" ------------------------------
# Date (copy formatting/value from B36 which already holds the same date).
$ws.Range("B36").Copy($ws.Range("B37"))

# Problem statement (new shared string -> becomes index 64).
$ws.Range("C36").Copy($ws.Range("C37"))
$ws.Range("C37").Value = "Maximum circular subarray sum"

# Solution complexity note (new shared string -> becomes index 65).
$ws.Range("F37").Value = "O(n^2)"

# "who / took help" column - reuses existing shared string "solved".
$ws.Range("G37").Value = "solved"

# --- Row 38: "Remove nth node from end of list." ---------------------------
$ws.Rows("38").RowHeight = 60

# Date (copy formatting from B36, then set the actual date value for row 38).
$ws.Range("B36").Copy($ws.Range("B38"))
$ws.Range("B38").Value = 45653

# Solution complexity note (new shared string -> becomes index 66).
# Written before C38 so the shared-string order matches the source edit.
$ws.Range("F38").Value = "O(n)"

# Problem statement with a bold "Leetcode" suffix (new shared string -> index 67).
$ws.Range("C35").Copy($ws.Range("C38"))
$ws.Range("C38").Value = "Remove nth node from end of list. Leetcode"
$ws.Range("C38").Characters(35, 8).Font.Bold = $true

# "who / took help" column - reuses existing shared string "solved and submitted".
$ws.Range("G38").Value = "solved and submitted"

# --- View state: scrolled so row 34 is at top, active cell on C42 ---------
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C42").Select() | Out-Null
